$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the "Bitbop DE" header block (old row 8 and its trailing spacer row 9)
# down by one row, opening up rows 6-8 under the "BitBop US" header for the
# new admin-link rows.
$ws.Rows("8").Insert()

# Column B needs to be wide enough to show the admin URLs.
$ws.Columns("B").ColumnWidth = 53

# --- BitBop US admin links (rows 6-7, row 8 left blank but styled) ---
$ws.Range("B6").Value = "http://www.stage-bitbop.com/admin"
$ws.Range("B7").Value = "http://www.bitbop.com/admin"
$ws.Range("B11").Value = "http://www.-bitbop.de/admin"
$ws.Range("B10").Value = "http://www.stage-bitbop.de/admin"

$ws.Range("A6").Value = "Staging admin"
$ws.Range("A7").Value = "Production admin"
$ws.Range("A10").Value = "Staging admin"
$ws.Range("A11").Value = "Production admin"

$ws.Hyperlinks.Add($ws.Range("B6"), "http://www.stage-bitbop.com/admin")
$ws.Hyperlinks.Add($ws.Range("B7"), "http://www.bitbop.com/admin")
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.stage-bitbop.de/admin")
$ws.Hyperlinks.Add($ws.Range("B11"), "http://www.-bitbop.de/admin")

# Blank spacer cell under the BitBop US links, carrying the hyperlink style
# like its neighbours.
$ws.Range("B8").Style = "Hyperlink"

# Leave the selection where the author ended up after editing.
[void]$ws.Range("A16").Select()
